$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packing Slip")

# New "COST" header in G10
$ws.Range("G10").Value = "COST"

# New cost formula in G11 = QUANTITY * PRICE
$ws.Range("G11").Formula = "=F11*E11"

# Move selection to the new cost cell
$ws.Activate()
$ws.Range("G11").Select()
